# Wrap every existing top-level shape on slide 1 into a single new group
# shape named "Group 53" (so the grouping operation mirrors the author's
# edit that combined all the "neighbors" help-slide shapes into one group).
#
# PowerPoint's shape-id allocator hands out the smallest unused id (>= 2)
# that isn't already present anywhere on the slide. Grouping the 22
# pre-existing shapes (ids 4,5,7,9,10,15,16,17,27,28,29,30,31,33,35,36,
# 39,40,45,49,50,52) would normally land on id 2 (the first free slot).
# To reproduce the target id (54) we first "burn" the 30 free ids that
# sort before 54 using throw-away textboxes, run the real Group() call
# (which then consumes id 54), and finally delete the throw-away shapes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Capture the names of all shapes currently on the slide - these are the
# ones that need to end up inside the new group.
$names = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $names += $s.Shapes.Item($i).Name
}

# Burn the 30 ids that precede 54 in the allocator's free-id sequence.
$burned = @()
for ($i = 1; $i -le 30; $i++) {
    $burned += $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
}

# Group all the original shapes together; this now receives id 54.
$range = $s.Shapes.Range($names)
$grp = $range.Group()
$grp.Name = "Group 53"

# Remove the throw-away shapes used purely to advance the id allocator.
foreach ($b in $burned) {
    $b.Delete()
}

Write-Output "Group created: id=$($grp.Id) name=$($grp.Name) shapes=$($grp.GroupItems.Count)"
